$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, copying the formatting (bold,
# centered, bordered style) used by the other header cells, e.g. G1 ("sum").
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the Save values for the data rows (plain numeric cells, no style).
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
